# Update "想去人数" (column F) figures to the freshly scraped counts.
# (gh-pages data refresh -- output generated at 456a3b4)

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value  = 248
$ws.Cells.Item(6, 6).Value  = 113
$ws.Cells.Item(9, 6).Value  = 4715
$ws.Cells.Item(10, 6).Value = 4715
$ws.Cells.Item(12, 6).Value = 449
$ws.Cells.Item(14, 6).Value = 613
$ws.Cells.Item(15, 6).Value = 4265
$ws.Cells.Item(17, 6).Value = 167
$ws.Cells.Item(20, 6).Value = 3471
$ws.Cells.Item(24, 6).Value = 3065
$ws.Cells.Item(25, 6).Value = 130
$ws.Cells.Item(29, 6).Value = 193
$ws.Cells.Item(32, 6).Value = 57
$ws.Cells.Item(36, 6).Value = 5484
$ws.Cells.Item(37, 6).Value = 813
$ws.Cells.Item(38, 6).Value = 396
$ws.Cells.Item(41, 6).Value = 47
$ws.Cells.Item(42, 6).Value = 1102
$ws.Cells.Item(46, 6).Value = 296
$ws.Cells.Item(47, 6).Value = 62
$ws.Cells.Item(48, 6).Value = 697

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value  = 11
$ws.Cells.Item(6, 6).Value  = 82
$ws.Cells.Item(8, 6).Value  = 39
$ws.Cells.Item(22, 6).Value = 736

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 199

# ---- 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 199
$ws.Cells.Item(4, 6).Value  = 248
$ws.Cells.Item(6, 6).Value  = 82
$ws.Cells.Item(7, 6).Value  = 113
$ws.Cells.Item(10, 6).Value = 4715
$ws.Cells.Item(11, 6).Value = 4715
$ws.Cells.Item(12, 6).Value = 39
$ws.Cells.Item(16, 6).Value = 449
$ws.Cells.Item(18, 6).Value = 613
$ws.Cells.Item(19, 6).Value = 4265
$ws.Cells.Item(21, 6).Value = 167
$ws.Cells.Item(23, 6).Value = 3471
$ws.Cells.Item(24, 6).Value = 3065
$ws.Cells.Item(25, 6).Value = 130
$ws.Cells.Item(28, 6).Value = 193
$ws.Cells.Item(36, 6).Value = 5484
$ws.Cells.Item(38, 6).Value = 813
$ws.Cells.Item(39, 6).Value = 396
$ws.Cells.Item(44, 6).Value = 47
$ws.Cells.Item(45, 6).Value = 1102
$ws.Cells.Item(48, 6).Value = 296
$ws.Cells.Item(49, 6).Value = 697
